$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores every "Price" cell as literal text (t="inlineStr"),
# even when the text looks like a plain number ("41.08", "22.30", ...). Excel.
# COM.Range.Value coerces a plain-numeric-looking string into a Double, which would
# flip these cells to numeric storage (and silently drop significant trailing zeros,
# e.g. "22.30" -> 22.3). Force Text format first so the literal string is kept,
# then restore the default "Normal" cell style so no stray number format lingers.
$textForceCells = "D5", "D6", "D9", "D10", "D11", "D14", "D16", "D20", "D22", "D24", "D25", "D29", "D31", "D32", "D33", "D35", "D36", "D37", "D39", "D42", "D44", "D47", "D48", "D49", "D51"
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "45.182.24"
$ws.Range("E2").Value = "  +5.76%  "

# Row 3
$ws.Range("D3").Value = "2.357.33"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "309.56"
$ws.Range("E5").Value = "  -0.70%  "

# Row 6
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "108.03"
$ws.Range("E6").Value = "  +1.20%  "

# Row 7
$ws.Range("E7").Value = "  +0.82%  "

# Row 8
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  +1.30%  "

# Row 10
$ws.Range("D10").Value = "41.08"
$ws.Range("E10").Value = "  +2.58%  "

# Row 11
$ws.Range("D11").Value = "0.0916"
$ws.Range("E11").Value = "  +0.19%  "

# Row 12
$ws.Range("E12").Value = "  +0.83%  "

# Row 13
$ws.Range("E13").Value = "  +2.14%  "

# Row 14
$ws.Range("D14").Value = "0.982"
$ws.Range("E14").Value = "  -0.64%  "

# Row 15
$ws.Range("D15").Value = "2.716.50"
$ws.Range("E15").Value = "  +1.93%  "

# Row 16
$ws.Range("D16").Value = "15.31"
$ws.Range("E16").Value = "  +0.00%  "

# Row 17
$ws.Range("D17").Value = "2.360.13"
$ws.Range("E17").Value = "  +1.74%  "

# Row 18
$ws.Range("D18").Value = "45.156.47"
$ws.Range("E18").Value = "  +5.74%  "

# Row 19
$ws.Range("E19").Value = "  -1.84%  "

# Row 20
$ws.Range("D20").Value = "13.97"
$ws.Range("E20").Value = "  +6.98%  "

# Row 21
$ws.Range("E21").Value = "  +0.63%  "

# Row 22
$ws.Range("D22").Value = "73.17"
$ws.Range("E22").Value = "  -0.56%  "

# Row 23
$ws.Range("E23").Value = "  +0.15%  "

# Row 24
$ws.Range("D24").Value = "259.35"
$ws.Range("E24").Value = "  -2.37%  "

# Row 25
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +3.34%  "

# Row 26
$ws.Range("E26").Value = "  -0.36%  "

# Row 27
$ws.Range("E27").Value = "  +0.73%  "

# Row 28
$ws.Range("E28").Value = "  -4.52%  "

# Row 29
$ws.Range("D29").Value = "2.35"
$ws.Range("E29").Value = "  +2.42%  "

# Row 30
$ws.Range("E30").Value = "  +10.27%  "

# Row 31
$ws.Range("D31").Value = "22.30"
$ws.Range("E31").Value = "  -0.66%  "

# Row 32
$ws.Range("D32").Value = "37.77"
$ws.Range("E32").Value = "  -2.69%  "

# Row 33
$ws.Range("D33").Value = "168.90"
$ws.Range("E33").Value = "  +1.57%  "

# Row 34
$ws.Range("E34").Value = "  +6.40%  "

# Row 35
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  +4.72%  "

# Row 37
$ws.Range("D37").Value = "4.80"
$ws.Range("E37").Value = "  +2.27%  "

# Row 38
$ws.Range("E38").Value = "  +5.20%  "

# Row 39
$ws.Range("D39").Value = "3.91"
$ws.Range("E39").Value = "  +6.52%  "

# Row 40
$ws.Range("E40").Value = "  -0.74%  "

# Row 41
$ws.Range("E41").Value = "  +7.56%  "

# Row 42
$ws.Range("D42").Value = "99.41"
$ws.Range("E42").Value = "  -4.26%  "

# Row 43
$ws.Range("E43").Value = "  -0.50%  "

# Row 44
$ws.Range("D44").Value = "69.59"
$ws.Range("E44").Value = "  -1.80%  "

# Row 45
$ws.Range("E45").Value = "  -0.76%  "

# Row 46
$ws.Range("E46").Value = "  +0.24%  "

# Row 47
$ws.Range("D47").Value = "81.98"
$ws.Range("E47").Value = "  +6.56%  "

# Row 48
$ws.Range("D48").Value = "112.11"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("D49").Value = "5.50"
$ws.Range("E49").Value = "  +5.11%  "

# Row 50
$ws.Range("D50").Value = "1.678.35"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51
$ws.Range("D51").Value = "9.13"
$ws.Range("E51").Value = "  +3.81%  "

# Restore default style on the text-forced cells (removes the temporary "@" format)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
